$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 20:35"

$ws.Range("B4").Value = 1395027
$ws.Range("C4").Value = 9193
$ws.Range("E4").Value = 1037571

$ws.Range("B9").Value = 178225
$ws.Range("C9").Value = 802
$ws.Range("D9").Value = 57785
$ws.Range("E9").Value = 93449
$ws.Range("F9").Value = 2542
$ws.Range("G9").Value = 348
$ws.Range("H9").Value = 26991

$ws.Range("F10").Value = 1539

$ws.Range("B59").Value = 5522
$ws.Range("C59").Value = 286
$ws.Range("D59").Value = 2192
$ws.Range("E59").Value = 3321

$ws.Range("B104").Value = 889
$ws.Range("C104").Value = 26
$ws.Range("E104").Value = 514
